$d = $word.ActiveDocument

# 1. Replace the "By: Rafil Yashooa" paragraph text with the new financial
#    status sentence. Find/Replace collapses the run set into a single run
#    and preserves the trailing bookmark (_GoBack) untouched.
$d.Content.Find.Execute(
    "By: Rafil Yashooa", $true, $false, $false, $false, $false, $true, 1,
    $false,
    "As for the financial status, everything is still the same. And no problems so far.",
    2) | Out-Null

# 2. Find that paragraph (now the last paragraph in the body) and split a
#    new empty paragraph after it.
$statusPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$statusPara.Range.InsertParagraphAfter() | Out-Null

# 3. The freshly created empty paragraph becomes the insertion target for
#    the restored "By: Rafil Yashooa" signature block, rebuilt with the
#    original run/proofErr structure (now split into two spell-checked
#    words, "Rafil" and "Yashooa", each wrapped in its own proofErr pair).
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$signatureXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr>
          </w:p>
          <w:p>
            <w:pPr><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr>
            <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve">By: </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>Rafil</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>Yashooa</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($signatureXml) | Out-Null

# 4. InsertXML cannot remove the document's final paragraph mark (Word never
#    allows the very last pilcrow of the body to be deleted), so it leaves a
#    blank paragraph behind at the end. Merge it away by deleting the break
#    between the new signature paragraph and that trailing leftover mark.
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
$last = $d.Paragraphs.Item($count)
$d.Range($secondToLast.Range.End - 1, $last.Range.End).Delete() | Out-Null
